# Fix two data-entry typos in the staff roster sheet:
#  - A2: "Mr.  BASKARAN V" (double space) -> "Mr. BASKARAN V" (single space)
#  - C4: "\/static/images/profile_photos/010/VEC-010-04-151.webp"
#        -> "/static/images/profile_photos/010/VEC-010-04-151.webp" (stray leading backslash removed)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C4").Value = "/static/images/profile_photos/010/VEC-010-04-151.webp"
$ws.Range("A2").Value = "Mr. BASKARAN V"
